# Updated cryptos list (price / 1h-volume refresh) with GitHub Actions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that looks like a plain number (e.g. "5.20", "42.29")
# while keeping it stored as text, matching the source data's inlineStr cells
# (Excel's normal .Value assignment would otherwise auto-convert these to
# numbers and drop the significant trailing digits). The NumberFormat/Style
# dance forces text storage without leaving a residual style on the cell.
function Set-TextCell($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '62.337.94'
$ws.Range("E2").Value = '  -2.14%  '
# Row 3
$ws.Range("D3").Value = '2.425.95'
$ws.Range("E3").Value = '  -2.17%  '
# Row 4
Set-TextCell "D4" '0.999'
$ws.Range("E4").Value = '  -0.06%  '
# Row 5
Set-TextCell "D5" '578.27'
$ws.Range("E5").Value = '  +0.42%  '
# Row 6
Set-TextCell "D6" '143.00'
$ws.Range("E6").Value = '  -4.00%  '
# Row 7
$ws.Range("E7").Value = '  +0.05%  '
# Row 8
$ws.Range("E8").Value = '  -2.82%  '
# Row 9
$ws.Range("D9").Value = '2.422.42'
$ws.Range("E9").Value = '  -2.22%  '
# Row 10
$ws.Range("E10").Value = '  -4.92%  '
# Row 11
$ws.Range("E11").Value = '  +0.81%  '
# Row 12
Set-TextCell "D12" '5.20'
$ws.Range("E12").Value = '  -1.43%  '
# Row 13
Set-TextCell "D13" '0.346'
$ws.Range("E13").Value = '  -3.56%  '
# Row 14
Set-TextCell "D14" '26.36'
$ws.Range("E14").Value = '  -3.09%  '
# Row 15
$ws.Range("E15").Value = '  -4.60%  '
# Row 16
$ws.Range("D16").Value = '2.869.86'
$ws.Range("E16").Value = '  -2.13%  '
# Row 17
$ws.Range("D17").Value = '62.178.77'
$ws.Range("E17").Value = '  -2.16%  '
# Row 18
$ws.Range("D18").Value = '2.417.78'
$ws.Range("E18").Value = '  -2.25%  '
# Row 19
Set-TextCell "D19" '10.96'
$ws.Range("E19").Value = '  -4.93%  '
# Row 20
Set-TextCell "D20" '7.08'
$ws.Range("E20").Value = '  -3.78%  '
# Row 21
Set-TextCell "D21" '329.46'
# Row 22
Set-TextCell "D22" '4.12'
$ws.Range("E22").Value = '  -1.99%  '
# Row 23
$ws.Range("E23").Value = '  -4.88%  '
# Row 24
$ws.Range("E24").Value = '  +0.26%  '
# Row 25
Set-TextCell "D25" '65.48'
$ws.Range("E25").Value = '  -0.98%  '
# Row 26
Set-TextCell "D26" '633.55'
$ws.Range("E26").Value = '  +0.94%  '
# Row 27
Set-TextCell "D27" '9.14'
$ws.Range("E27").Value = '  +6.51%  '
# Row 28
$ws.Range("D28").Value = '2.540.85'
$ws.Range("E28").Value = '  -2.61%  '
# Row 29
$ws.Range("E29").Value = '  +0.03%  '
# Row 30
$ws.Range("D30").Value = '0.0₃0948'
$ws.Range("E30").Value = '  -9.00%  '
# Row 31
$ws.Range("E31").Value = '  -6.86%  '
# Row 32
Set-TextCell "D32" '8.03'
$ws.Range("E32").Value = '  -4.21%  '
# Row 33
$ws.Range("E33").Value = '  -1.45%  '
# Row 34
Set-TextCell "D34" '0.139'
$ws.Range("E34").Value = '  -3.26%  '
# Row 35
Set-TextCell "D35" '4.95'
$ws.Range("E35").Value = '  -5.46%  '
# Row 36
$ws.Range("E36").Value = '  +0.24%  '
# Row 37
Set-TextCell "D37" '1.45'
$ws.Range("E37").Value = '  -4.97%  '
# Row 38
Set-TextCell "D38" '0.374'
$ws.Range("E38").Value = '  -2.71%  '
# Row 39
Set-TextCell "D39" '148.19'
$ws.Range("E39").Value = '  +0.26%  '
# Row 40
$ws.Range("E40").Value = '  -2.61%  '
# Row 41
Set-TextCell "D41" '5.26'
$ws.Range("E41").Value = '  -3.74%  '
# Row 42
$ws.Range("E42").Value = '  -4.25%  '
# Row 43
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell "D43" '42.29'
$ws.Range("E43").Value = '  +1.00%  '
# Row 44
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextCell "D44" '0.999'
$ws.Range("E44").Value = '  +0.01%  '
# Row 45
Set-TextCell "D45" '2.48'
$ws.Range("E45").Value = '  -9.11%  '
# Row 46
Set-TextCell "D46" '143.30'
$ws.Range("E46").Value = '  -4.35%  '
# Row 47
$ws.Range("E47").Value = '  -2.91%  '
# Row 48
Set-TextCell "D48" '0.0519'
$ws.Range("E48").Value = '  -4.43%  '
# Row 49
Set-TextCell "D49" '0.594'
$ws.Range("E49").Value = '  -1.98%  '
# Row 50
Set-TextCell "D50" '19.45'
$ws.Range("E50").Value = '  -8.74%  '
# Row 51
$ws.Range("D51").Value = '0.0₆0239'
$ws.Range("E51").Value = '  +9.04%  '
